# Update the "想去人数" (want-to-go count) column F values across the
# workbook's sheets to match the newly scraped figures.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1298
$ws.Range("F4").Value = 944
$ws.Range("F5").Value = 982
$ws.Range("F6").Value = 1754
$ws.Range("F7").Value = 392
$ws.Range("F8").Value = 1171
$ws.Range("F9").Value = 52
$ws.Range("F10").Value = 9
$ws.Range("F11").Value = 119
$ws.Range("F12").Value = 269
$ws.Range("F13").Value = 54
$ws.Range("F14").Value = 81
$ws.Range("F15").Value = 653
$ws.Range("F16").Value = 143
$ws.Range("F20").Value = 323
$ws.Range("F21").Value = 118
$ws.Range("F22").Value = 651
$ws.Range("F23").Value = 19
$ws.Range("F24").Value = 634
$ws.Range("F25").Value = 144
$ws.Range("F27").Value = 853
$ws.Range("F28").Value = 304
$ws.Range("F29").Value = 135
$ws.Range("F30").Value = 32
$ws.Range("F31").Value = 256
$ws.Range("F32").Value = 9
$ws.Range("F34").Value = 403

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 314
$ws.Range("F11").Value = 117

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 302

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 302
$ws.Range("F4").Value = 1298
$ws.Range("F5").Value = 944
$ws.Range("F6").Value = 982
$ws.Range("F7").Value = 1754
$ws.Range("F8").Value = 392
$ws.Range("F9").Value = 1171
$ws.Range("F10").Value = 52
$ws.Range("F12").Value = 9
$ws.Range("F13").Value = 119
$ws.Range("F14").Value = 269
$ws.Range("F15").Value = 54
$ws.Range("F16").Value = 81
$ws.Range("F17").Value = 653
$ws.Range("F18").Value = 143
$ws.Range("F22").Value = 314
$ws.Range("F25").Value = 323
$ws.Range("F29").Value = 118
$ws.Range("F30").Value = 651
$ws.Range("F31").Value = 19
$ws.Range("F32").Value = 634
$ws.Range("F33").Value = 144
$ws.Range("F35").Value = 853
$ws.Range("F36").Value = 304
$ws.Range("F39").Value = 135
$ws.Range("F40").Value = 32
$ws.Range("F41").Value = 256
$ws.Range("F43").Value = 117
$ws.Range("F44").Value = 117
$ws.Range("F45").Value = 9
$ws.Range("F48").Value = 403
